# Updates cryptos list: refreshed prices/volumes for most rows, plus a
# re-ranking of several coins (rows 31-34 and 45-46 swap place/identity).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.777.74'
$ws.Range("E2").Value = '  +1.47%  '
$ws.Range("D3").Value = '3.364.16'
$ws.Range("E3").Value = '  +8.00%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.37'
$ws.Range("E5").Value = '  +8.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.11'
$ws.Range("E6").Value = '  +3.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.40'
$ws.Range("E7").Value = '  +26.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.392'
$ws.Range("E8").Value = '  +2.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.875'
$ws.Range("E10").Value = '  +11.36%  '
$ws.Range("D11").Value = '3.359.84'
$ws.Range("E11").Value = '  +7.94%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").Value = '98.548.11'
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.95'
$ws.Range("E14").Value = '  +6.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000248'
$ws.Range("E15").Value = '  +2.78%  '
$ws.Range("D16").Value = '3.984.29'
$ws.Range("E16").Value = '  +7.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.50'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("D18").Value = '3.366.08'
$ws.Range("E18").Value = '  +7.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.57'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.20'
$ws.Range("E20").Value = '  +4.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '492.47'
$ws.Range("E21").Value = '  -7.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.20'
$ws.Range("E22").Value = '  +9.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000212'
$ws.Range("E23").Value = '  +9.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.46'
$ws.Range("E24").Value = '  +7.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.68'
$ws.Range("E25").Value = '  +3.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.68'
$ws.Range("E26").Value = '  +0.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.93'
$ws.Range("E27").Value = '  +3.03%  '
$ws.Range("D28").Value = '3.538.10'
$ws.Range("E28").Value = '  +7.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.285'
$ws.Range("E29").Value = '  +20.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.137'
$ws.Range("E31").Value = '  +10.17%  '
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.187'
$ws.Range("E32").Value = '  +7.51%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +15.08%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.55'
$ws.Range("E34").Value = '  +6.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.80'
$ws.Range("E35").Value = '  +4.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.151'
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.30'
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.97'
$ws.Range("E38").Value = '  +5.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.463'
$ws.Range("E39").Value = '  +6.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '497.61'
$ws.Range("E40").Value = '  +3.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.84'
$ws.Range("E41").Value = '  +2.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.81'
$ws.Range("E42").Value = '  +6.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.27'
$ws.Range("E43").Value = '  +3.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  +4.27%  '
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.780'
$ws.Range("E46").Value = '  +13.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.04'
$ws.Range("E47").Value = '  -1.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.94'
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.831'
$ws.Range("E49").Value = '  +14.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.65'
$ws.Range("E50").Value = '  +4.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.07'
$ws.Range("E51").Value = '  +3.61%  '
